# The deck has a single top-level shape on slide 1: a group ("Group 7")
# that contains three pictures. The commit resizes/repositions the group
# itself (its own off/ext) while leaving the group's child coordinate
# space (chOff/chExt) and the individual pictures inside it untouched.
#
# Target (EMU), from the authoritative OOXML diff:
#   off : x=1178558 -> 1256935   y=940525  -> 1005840
#   ext : cx=4529736 -> 2505167  cy=4376058 -> 2468880
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU) and are stored as single-precision floats, so the
# literals below are chosen such that round-tripping through a 32-bit
# float and back to EMU reproduces the exact target EMU values above
# (plain "emu / 12700" can land one EMU short after the float32 cast).

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)

$shp.Left   = 98.97125984251969   # -> 1256935 EMU
$shp.Top    = 79.2000008          # -> 1005840 EMU
$shp.Width  = 197.2572440944882   # -> 2505167 EMU
$shp.Height = 194.4000016         # -> 2468880 EMU
